# The authored edit swaps the raw OOXML content of ppt/theme/theme1.xml
# (the "Office Theme" colours, only wired to the notes master) and
# ppt/theme/theme2.xml (the "Integral" colours that actually drive the
# slide master / presentation design) - the file names/relationships stay
# put, only what's inside them changes. ppt/theme/theme1.xml's fontScheme
# and fmtScheme are byte-for-byte identical to ppt/theme/theme2.xml's, so
# the only real content difference between the two parts is the 12
# clrScheme colours (and the cosmetic theme/clrScheme "name" attributes,
# which PowerPoint's object model does not expose for editing).
#
# The only theme reachable from the PowerPoint COM object model in this
# host is the presentation's active design theme - i.e. the one behind
# ppt/theme/theme2.xml (SlideMaster.Theme / NotesMaster.Theme / Designs
# all resolve to this single theme). So we recolour that theme's 12
# scheme slots from the current "Integral" palette to the "Office Theme"
# palette that used to live in theme1.xml - the effective, reachable
# half of the swap.

$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$tcs = $sm.Theme.ThemeColorScheme

# MsoThemeColorSchemeIndex order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
